$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodo Mora (column E) values for rows 16-29, now listed in ascending
# chronological order (1702 .. 1712, 1801, 1802, 1803) instead of the
# previous descending order.
$periodos = @("1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712","1801","1802","1803")

$row = 16
foreach ($periodo in $periodos) {
    $ws.Range("E$row").Value = $periodo
    $row = $row + 1
}

# Salario Basico (column G) updated from 1092000 to 848000 for rows 16-29.
for ($row = 16; $row -le 29; $row++) {
    $ws.Range("G$row").Value = 848000
}
